# Model - need fix lasso
# Insert a new "Food" industry row after the header row, shifting all
# subsequent industry rows down by one, and append a new "Other" row
# at the bottom (mirroring the previously-last row's C/D/E stats),
# then refresh the simulated "Ann mean"/"Ann Sharpe" figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (pushes existing data rows 2-30 down to 3-31)
$ws.Rows(2).Insert()

# New "Food" row
$ws.Range("A2").Value = "Food "
$ws.Range("B2").Value = 8.28858394160584
$ws.Range("C2").Value = 15.02256728639485
$ws.Range("D2").Value = -18.15
$ws.Range("E2").Value = 19.89
$ws.Range("F2").Value = 0.5517421745291416

# Refreshed "Ann mean" (column B) and "Ann Sharpe" (column F) values for the
# rows that shifted down (their Ann vol/Minimum/Maximum stats - C, D, E -
# are carried over unchanged from the row above's previous position).
$ws.Range("B3").Value = 8.527357664233577
$ws.Range("F3").Value = 0.4839550027404988

$ws.Range("B4").Value = 11.78785401459854
$ws.Range("F4").Value = 0.5617917542693615

$ws.Range("B5").Value = 8.420496350364957
$ws.Range("F5").Value = 0.3387539152170506

$ws.Range("B6").Value = 6.339328467153287
$ws.Range("F6").Value = 0.3152426804573548

$ws.Range("B7").Value = 6.650277372262765
$ws.Range("F7").Value = 0.4036183776638168

$ws.Range("B8").Value = 8.033518248175184
$ws.Range("F8").Value = 0.3634139689040561

$ws.Range("B9").Value = 7.810861313868612
$ws.Range("F9").Value = 0.4578766144992515

$ws.Range("B10").Value = 6.23737226277372
$ws.Range("F10").Value = 0.3265194484806868

$ws.Range("B11").Value = 8.00899270072992
$ws.Range("F11").Value = 0.3294652170272067

$ws.Range("B12").Value = 6.281518248175189
$ws.Range("F12").Value = 0.3028769948412639

$ws.Range("B13").Value = 3.543591240875913
$ws.Range("F13").Value = 0.1403166371945758

$ws.Range("B14").Value = 6.75100729927007
$ws.Range("F14").Value = 0.3189773036310222

$ws.Range("B15").Value = 8.450627737226274
$ws.Range("F15").Value = 0.3930230334342855

$ws.Range("B16").Value = 5.447824817518251
$ws.Range("F16").Value = 0.2352790498340036

$ws.Range("B17").Value = 8.624759124087598
$ws.Range("F17").Value = 0.3959365397008792

$ws.Range("B18").Value = 6.583883211678836
$ws.Range("F18").Value = 0.2552022368620829

$ws.Range("B19").Value = 8.826569343065701
$ws.Range("F19").Value = 0.250149423045686

$ws.Range("B20").Value = 7.811737226277371
$ws.Range("F20").Value = 0.421940637358976

$ws.Range("B21").Value = 5.797489051094894
$ws.Range("F21").Value = 0.4197894450958571

$ws.Range("B22").Value = 6.260846715328465
$ws.Range("F22").Value = 0.390715583314208

$ws.Range("B23").Value = 8.200992700729923
$ws.Range("F23").Value = 0.3630200451129147

$ws.Range("B24").Value = 6.965781021897816
$ws.Range("F24").Value = 0.2986415029873843

$ws.Range("B25").Value = 6.120525547445252
$ws.Range("F25").Value = 0.3497799068549914

$ws.Range("B26").Value = 7.006598540145987
$ws.Range("F26").Value = 0.3531135399917951

$ws.Range("B27").Value = 7.471883211678826
$ws.Range("F27").Value = 0.3850874385277825

$ws.Range("B28").Value = 7.972554744525541
$ws.Range("F28").Value = 0.4305174614045045

$ws.Range("B29").Value = 8.467795620437951
$ws.Range("F29").Value = 0.4006927836192365

$ws.Range("B30").Value = 7.310715328467154
$ws.Range("F30").Value = 0.3902545721921423

# New trailing "Other" row (row 31), carrying the C/D/E stats that used to
# belong to the old last row ("Other" at row 30) forward, with refreshed
# Ann mean / Ann Sharpe values.
$ws.Range("A31").Value = "Other"
$ws.Range("B31").Value = 4.529167883211674
$ws.Range("C31").Value = 20.15123387703878
$ws.Range("D31").Value = -28.09
$ws.Range("E31").Value = 19.96
$ws.Range("F31").Value = 0.2247588366473386
